$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 924.25
$ws.Range("I2").Value = 932.3333
$ws.Range("K2").Value = 932.3333
$ws.Range("M2").Value = -819.3333

# Row 6
$ws.Range("H6").Value = 899
$ws.Range("I6").Value = 899
$ws.Range("K6").Value = 2697
$ws.Range("M6").Value = -2585

# Row 51
$ws.Range("H51").Value = 16750
$ws.Range("I51").Value = 13000
$ws.Range("K51").Value = 13000
$ws.Range("M51").Value = -12516

# Row 74
$ws.Range("H74").Value = 8333
$ws.Range("I74").Value = 8333
$ws.Range("K74").Value = 8333
$ws.Range("M74").Value = -7397

# Row 77
$ws.Range("H77").Value = 8333
$ws.Range("I77").Value = 8333
$ws.Range("K77").Value = 41665
$ws.Range("M77").Value = -36985

# Row 100
$ws.Range("H100").Value = 3070.4285
$ws.Range("I100").Value = 2299.6
$ws.Range("K100").Value = 2299.6
$ws.Range("M100").Value = -1758.6

# Row 111
$ws.Range("H111").Value = 2925.2222
$ws.Range("I111").Value = 4065.4
$ws.Range("K111").Value = 12196.2
$ws.Range("M111").Value = -9129.200000000001

# Row 116
$ws.Range("H116").Value = 4666.3335
$ws.Range("J116").Value = 5499.5
$ws.Range("L116").Value = 5499.5
$ws.Range("N116").Value = -12383.5

# Row 137
$ws.Range("H137").Value = 2837.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2837.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 8511.75
$ws.Range("N137").Value = -13611.75
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 159.71428
$ws.Range("I5").Value = 215.6
$ws.Range("K5").Value = 215.6
$ws.Range("M5").Value = -103.6

# Row 63
$ws.Range("H63").Value = 5142.143
$ws.Range("I63").Value = 1498.75
$ws.Range("K63").Value = 1498.75
$ws.Range("M63").Value = -812.75

# Row 66
$ws.Range("H66").Value = 5142.143
$ws.Range("I66").Value = 1498.75
$ws.Range("K66").Value = 7493.75
$ws.Range("M66").Value = -4061.75

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 1499.5714
$ws.Range("I29").Value = 1499.5714
$ws.Range("K29").Value = 1499.5714
$ws.Range("M29").Value = -1210.5714

# Row 36
$ws.Range("H36").Value = 9647.1
$ws.Range("I36").Value = 7184
$ws.Range("K36").Value = 7184
$ws.Range("M36").Value = -6650

# Row 81
$ws.Range("H81").Value = 31899.6
$ws.Range("J81").Value = 27374.5
$ws.Range("L81").Value = 27374.5
$ws.Range("N81").Value = -29496.5

# Row 84
$ws.Range("H84").Value = 31899.6
$ws.Range("J84").Value = 27374.5
$ws.Range("L84").Value = 82123.5
$ws.Range("N84").Value = -92731.5

# Row 99
$ws.Range("H99").Value = 1534.25
$ws.Range("I99").Value = 1534.25
$ws.Range("K99").Value = 1534.25
$ws.Range("M99").Value = -36.25

# Row 106
$ws.Range("H106").Value = 21910.143
$ws.Range("J106").Value = 21910.143
$ws.Range("L106").Value = 21910.143
$ws.Range("N106").Value = -24434.143

# Row 134
$ws.Range("H134").Value = 710
$ws.Range("I134").Value = 710
$ws.Range("K134").Value = 2130
$ws.Range("M134").Value = 405

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 586.3333
$ws.Range("I7").Value = 678.6
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 678.6
$ws.Range("L7").Value = 125
$ws.Range("M7").Value = -565.6
$ws.Range("N7").Value = -351

# Row 31
$ws.Range("H31").Value = 8127.5713
$ws.Range("I31").Value = 5012
$ws.Range("J31").Value = 8646.833000000001
$ws.Range("K31").Value = 5012
$ws.Range("L31").Value = 8646.833000000001
$ws.Range("M31").Value = -4717
$ws.Range("N31").Value = -9236.833000000001

# Row 34
$ws.Range("H34").Value = 8127.5713
$ws.Range("I34").Value = 5012
$ws.Range("J34").Value = 8646.833000000001
$ws.Range("K34").Value = 5012
$ws.Range("L34").Value = 8646.833000000001
$ws.Range("M34").Value = -4810
$ws.Range("N34").Value = -9050.833000000001

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 103
$ws.Range("H103").Value = 39973.25
$ws.Range("I103").Value = 39973.25
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 39973.25
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -38801.25
$ws.Range("N103").ClearContents()

# Row 134
$ws.Range("H134").Value = 4437.933
$ws.Range("I134").Value = 2350.818
$ws.Range("K134").Value = 7052.454000000001
$ws.Range("M134").Value = -4517.454000000001

# Row 141
$ws.Range("H141").Value = 515170.56
$ws.Range("J141").Value = 515170.56
$ws.Range("L141").Value = 515170.56
$ws.Range("N141").Value = -525530.5600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 4084.35
$ws.Range("J55").Value = 4142.559
$ws.Range("L55").Value = 12427.677
$ws.Range("N55").Value = -12781.677

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 134
$ws.Range("H134").Value = 987.4
$ws.Range("I134").Value = 987.4
$ws.Range("K134").Value = 2962.2
$ws.Range("M134").Value = 2107.8

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2670.5
$ws.Range("I97").Value = 2985.2856
$ws.Range("J97").Value = 2229.8
$ws.Range("K97").Value = 2985.2856
$ws.Range("L97").Value = 2229.8
$ws.Range("M97").Value = -2489.2856
$ws.Range("N97").Value = -3221.8

# Row 126
$ws.Range("H126").Value = 1786.6666
$ws.Range("I126").Value = 1761.2
$ws.Range("J126").Value = 1914
$ws.Range("K126").Value = 5283.6
$ws.Range("L126").Value = 5742
$ws.Range("M126").Value = -2813.6
$ws.Range("N126").Value = -10682

# Row 132
$ws.Range("H132").Value = 2265.375
$ws.Range("I132").Value = 1824.6
$ws.Range("K132").Value = 5473.799999999999
$ws.Range("M132").Value = -2943.799999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -105

# Row 27
$ws.Range("H27").Value = 550
$ws.Range("I27").Value = 400
$ws.Range("K27").Value = 400
$ws.Range("M27").Value = -293

# Row 30
$ws.Range("H30").Value = 1088.5
$ws.Range("I30").Value = 906.2
$ws.Range("K30").Value = 906.2
$ws.Range("M30").Value = -798.2

# Row 82
$ws.Range("H82").Value = 3666
$ws.Range("I82").Value = 3399.4
$ws.Range("K82").Value = 3399.4
$ws.Range("M82").Value = -3038.4

# Row 85
$ws.Range("H85").Value = 3666
$ws.Range("I85").Value = 3399.4
$ws.Range("K85").Value = 3399.4
$ws.Range("M85").Value = -2151.4

# Row 100
$ws.Range("H100").Value = 2006.9286
$ws.Range("I100").Value = 2172.9092
$ws.Range("J100").Value = 1398.3334
$ws.Range("K100").Value = 2172.9092
$ws.Range("L100").Value = 1398.3334
$ws.Range("M100").Value = -1631.9092
$ws.Range("N100").Value = -2480.3334

# Row 132
$ws.Range("H132").Value = 3607.8
$ws.Range("I132").Value = 3607.8
$ws.Range("K132").Value = 10823.4
$ws.Range("M132").Value = -8293.400000000001

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 122
$ws.Range("I122").Value = 1563.3636
$ws.Range("J122").Value = 2112
$ws.Range("K122").Value = 4690.0908
$ws.Range("L122").Value = 6336
$ws.Range("M122").Value = -2240.0908
$ws.Range("N122").Value = -11236

# Row 137
$ws.Range("H137").Value = 42715
$ws.Range("J137").Value = 42715
$ws.Range("L137").Value = 42715
$ws.Range("N137").Value = -52915
